# Applies the team-name corrections described by the commit:
#   - fixes the misspelled team "Spirt" -> "Spirit" on every sheet
#   - on the "传奇组" (Legend group) sheet only, renames "Tyloo" -> "二次元"
#     and "Astralls" -> "原神"
# and restores the selection/active-cell state left behind by the edit.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 挑战组 (Challenger group) -------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B3").Value = "Spirit"
$ws1.Range("B8").Value = "Spirit"
$ws1.Range("B10").Value = "Spirit"
$ws1.Range("C12").Value = "Spirit"
$ws1.Range("B15").Value = "Spirit"
$ws1.Range("B16").Value = "Spirit"
$ws1.Range("C17").Value = "Spirit"

# --- Sheet 2: 中坚组 (Core group) --------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B3").Value = "Spirit"
$ws2.Range("B8").Value = "Spirit"
$ws2.Range("B10").Value = "Spirit"
$ws2.Range("C12").Value = "Spirit"
$ws2.Range("B15").Value = "Spirit"
$ws2.Range("B16").Value = "Spirit"
$ws2.Range("C17").Value = "Spirit"

# --- Sheet 3: 传奇组 (Legend group) ------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "二次元"
$ws3.Range("C2").Value = "原神"
$ws3.Range("B3").Value = "Spirit"
$ws3.Range("B6").Value = "二次元"
$ws3.Range("B8").Value = "Spirit"
$ws3.Range("C8").Value = "原神"
$ws3.Range("B10").Value = "Spirit"
$ws3.Range("B12").Value = "二次元"
$ws3.Range("C12").Value = "Spirit"
$ws3.Range("B14").Value = "二次元"
$ws3.Range("B15").Value = "Spirit"
$ws3.Range("B16").Value = "Spirit"
$ws3.Range("B17").Value = "二次元"
$ws3.Range("C17").Value = "Spirit"

# --- Restore the selections left active in each sheet ------------------
$ws1.Activate() | Out-Null
$ws1.Range("B3").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("O12").Select() | Out-Null
